{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst newValues = [\n  [\"44+22=66\", \"30+58=88\", \"70-3=67\", \"57+42=99\", \"69-24=45\"],\n  [\"22-20=2\", \"42-11=31\", \"46-22=24\", \"70-39=31\", \"74-46=28\"],\n  [\"73-18=55\", \"75-6=69\", \"90-54=36\", \"55-19=36\", \"42+6=48\"],\n  [\"55-53=2\", \"74+9=83\", \"27+49=76\", \"85-0=85\", \"59+22=81\"],\n  [\"49+5=54\", \"68-44=24\", \"87-56=31\", \"68-17=51\", \"51+23=74\"],\n  [\"89-32=57\", \"75+10=85\", \"59-18=41\", \"76+20=96\", \"22+75=97\"],\n  [\"98-12=86\", \"96+1=97\", \"40+29=69\", \"45-1=44\", \"18+67=85\"],\n  [\"44-29=15\", \"31-12=19\", \"3+45=48\", \"80-68=12\", \"77-49=28\"],\n  [\"43-7=36\", \"99-78=21\", \"2-2=0\", \"86-11=75\", \"51-38=13\"],\n  [\"71-18=53\", \"44-5=39\", \"10+74=84\", \"84-31=53\", \"27+21=48\"],\n  [\"43+29=72\", \"65-52=13\", \"20+53=73\", \"98-8=90\", \"5-1=4\"],\n  [\"74-44=30\", \"39+37=76\", \"65-33=32\", \"15+23=38\", \"87-51=36\"],\n  [\"55-26=29\", \"63-22=41\", \"99-69=30\", \"45-15=30\", \"28+37=65\"],\n  [\"94-9=85\", \"82-70=12\", \"48+15=63\", \"61-28=33\", \"86-26=60\"],\n  [\"27+53=80\", \"9+67=76\", \"16+48=64\", \"0+10=10\", \"47+23=70\"],\n  [\"9+67=76\", \"60+15=75\", \"52+42=94\", \"41+12=53\", \"96-45=51\"],\n  [\"40-34=6\", \"9+29=38\", \"6+34=40\", \"23+38=61\", \"50-21=29\"],\n  [\"97-31=66\", \"33+17=50\", \"67+18=85\", \"48+46=94\", \"55-24=31\"],\n  [\"29+16=45\", \"65-48=17\", \"16+60=76\", \"17+78=95\", \"5+0=5\"],\n  [\"29-24=5\", \"33+64=97\", \"28+54=82\", \"35+17=52\", \"25+52=77\"]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"44+22=66\", \"30+58=88\", \"70-3=67\", \"57+42=99\", \"69-24=45\"),\n    @(\"22-20=2\", \"42-11=31\", \"46-22=24\", \"70-39=31\", \"74-46=28\"),\n    @(\"73-18=55\", \"75-6=69\", \"90-54=36\", \"55-19=36\", \"42+6=48\"),\n    @(\"55-53=2\", \"74+9=83\", \"27+49=76\", \"85-0=85\", \"59+22=81\"),\n    @(\"49+5=54\", \"68-44=24\", \"87-56=31\", \"68-17=51\", \"51+23=74\"),\n    @(\"89-32=57\", \"75+10=85\", \"59-18=41\", \"76+20=96\", \"22+75=97\"),\n    @(\"98-12=86\", \"96+1=97\", \"40+29=69\", \"45-1=44\", \"18+67=85\"),\n    @(\"44-29=15\", \"31-12=19\", \"3+45=48\", \"80-68=12\", \"77-49=28\"),\n    @(\"43-7=36\", \"99-78=21\", \"2-2=0\", \"86-11=75\", \"51-38=13\"),\n    @(\"71-18=53\", \"44-5=39\", \"10+74=84\", \"84-31=53\", \"27+21=48\"),\n    @(\"43+29=72\", \"65-52=13\", \"20+53=73\", \"98-8=90\", \"5-1=4\"),\n    @(\"74-44=30\", \"39+37=76\", \"65-33=32\", \"15+23=38\", \"87-51=36\"),\n    @(\"55-26=29\", \"63-22=41\", \"99-69=30\", \"45-15=30\", \"28+37=65\"),\n    @(\"94-9=85\", \"82-70=12\", \"48+15=63\", \"61-28=33\", \"86-26=60\"),\n    @(\"27+53=80\", \"9+67=76\", \"16+48=64\", \"0+10=10\", \"47+23=70\"),\n    @(\"9+67=76\", \"60+15=75\", \"52+42=94\", \"41+12=53\", \"96-45=51\"),\n    @(\"40-34=6\", \"9+29=38\", \"6+34=40\", \"23+38=61\", \"50-21=29\"),\n    @(\"97-31=66\", \"33+17=50\", \"67+18=85\", \"48+46=94\", \"55-24=31\"),\n    @(\"29+16=45\", \"65-48=17\", \"16+60=76\", \"17+78=95\", \"5+0=5\"),\n    @(\"29-24=5\", \"33+64=97\", \"28+54=82\", \"35+17=52\", \"25+52=77\"),\n)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
